# Updated symbol list (price/volume refresh + WRX/ONE swap) with GitHub Actions
# Values that look numeric are prefixed with a literal apostrophe so Excel
# stores them as text (matching the source sheet's text-typed Price/Hora columns)
# instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''242.94'
$ws.Range("G2").Value = '''20'
$ws.Range("G3").Value = '''20'
$ws.Range("D4").Value = '''5.245'
$ws.Range("G4").Value = '''20'
$ws.Range("G5").Value = '''20'
$ws.Range("D6").Value = '''3.371'
$ws.Range("G6").Value = '''20'
$ws.Range("D7").Value = '''6.375'
$ws.Range("G7").Value = '''20'
$ws.Range("D8").Value = '''0.8056'
$ws.Range("G8").Value = '''20'
$ws.Range("D9").Value = '''0.9397'
$ws.Range("G9").Value = '''20'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.0005803'
$ws.Range("E10").Value = '9OneONE'
$ws.Range("G10").Value = '''20'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1422'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("G11").Value = '''20'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07313'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("G12").Value = '''20'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '''0.03153'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G13").Value = '''20'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.03014'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("G14").Value = '''20'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '''0.09271'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("G15").Value = '''20'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '''3.620'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("G16").Value = '''20'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '''0.001657'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("G17").Value = '''20'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '''0.04713'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("G18").Value = '''20'
$ws.Range("D19").Value = '''0.006361'
$ws.Range("G19").Value = '''20'
$ws.Range("D20").Value = '''0.004972'
$ws.Range("G20").Value = '''20'
$ws.Range("D21").Value = '''0.001043'
$ws.Range("G21").Value = '''20'
$ws.Range("D22").Value = '''0.0001501'
$ws.Range("G22").Value = '''20'
$ws.Range("G23").Value = '''20'
$ws.Range("D24").Value = '''3.766'
$ws.Range("G24").Value = '''20'
$ws.Range("G25").Value = '''20'
$ws.Range("D26").Value = '''0.3228'
$ws.Range("G26").Value = '''20'
$ws.Range("G27").Value = '''20'
$ws.Range("G28").Value = '''20'
$ws.Range("G29").Value = '''20'
$ws.Range("G30").Value = '''20'
$ws.Range("G31").Value = '''20'
$ws.Range("G32").Value = '''20'
$ws.Range("G33").Value = '''20'
$ws.Range("G34").Value = '''20'
$ws.Range("G35").Value = '''20'
$ws.Range("G36").Value = '''20'
$ws.Range("G37").Value = '''20'
$ws.Range("G38").Value = '''20'
$ws.Range("G39").Value = '''20'
$ws.Range("D40").Value = '''0.03907'
$ws.Range("G40").Value = '''20'
$ws.Range("D41").Value = '''0.006882'
$ws.Range("G41").Value = '''20'
$ws.Range("G42").Value = '''20'
$ws.Range("D43").Value = '''0.1030'
$ws.Range("G43").Value = '''20'
$ws.Range("D44").Value = '''0.008261'
$ws.Range("G44").Value = '''20'
$ws.Range("D45").Value = '''0.00005941'
$ws.Range("G45").Value = '''20'
$ws.Range("G46").Value = '''20'
$ws.Range("G47").Value = '''20'
$ws.Range("D48").Value = '''0.6827'
$ws.Range("G48").Value = '''20'
$ws.Range("D49").Value = '''0.05965'
$ws.Range("G49").Value = '''20'
$ws.Range("G50").Value = '''20'
$ws.Range("G51").Value = '''20'
